$d = $word.ActiveDocument

# Remove the comma: "mind, but" -> "mind but"
$d.Content.Find.Execute("mind, but", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mind but", 2) | Out-Null

# Word's editor records the last edit position with the special "_GoBack"
# bookmark. Remove any existing one and re-create it at the location of
# the edit we just made (right after "mind but").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$findRange = $d.Content
$findRange.Find.Execute("mind but", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
